# Generate Report for Handback
# Update the "Status" for the 23d47892... file (row 3) to reflect a failed
# handback transform, and record the error detail for each locale sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns for zh-cn (B) and de-de (C) on row 3
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"

# zh-cn sheet: Status (C3) + Error Detail (K3) for row 3
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("K3").Value = "Handback file name: fyoiwcr4.sjl is different with handoff file name: 23d47892-f416-45d4-9a8b-decff2b66808.fe660805b67a5cbe8577ad3163b92ce04c062868.zh-cn."

# de-de sheet: Status (C3) + Error Detail (K3) for row 3
$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("K3").Value = "Handback file name: fyoiwcr4.sjl is different with handoff file name: 23d47892-f416-45d4-9a8b-decff2b66808.fe660805b67a5cbe8577ad3163b92ce04c062868.de-de."
